$d = $word.ActiveDocument

# --- Step 1: change the existing final paragraph's text ---
# The document currently ends with a paragraph: <w:tab/>"You can look at " + bookmark.
# We replace the sentence with "Items to touch:". This paragraph is kept as-is (still
# holding the leading tab); the new content is appended after it, and finally a fresh
# paragraph (tab only, bold) takes over hosting the "_GoBack" bookmark.
$null = $d.Content.Find.Execute("You can look at ", $true, $false, $false, $false, $false, $true, 1, $false, "Items to touch:", 2)

$global:base = $d.Paragraphs.Last

function New-Para([string]$text) {
    $r = $global:base.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $np = $d.Paragraphs.Last
    if ($text -ne "") {
        $nr = $np.Range
        $nr.Collapse(0)
        $nr.InsertBefore($text)
    }
    $global:base = $np
    return $np
}

# ===== Build all new paragraphs first (plain, in final order) =====
$pShoe      = New-Para("Shoe marks")
$pShoeDesc  = New-Para("“Hmm, big size. They are so wet and dirty. No wonder, because it was raining yesterday. If you look at them accurate you can notice that the toe of the boot hardly touched the floor. Also, steps are very gradual, killer was not worried at all. Strange, very strange”")
$pHarpoon   = New-Para("Harpoon")
$pHarpoonD  = New-Para("“What a strength is needed to nail a person with harpoon. It can’t be a woman.”")
$pHook      = New-Para("Message holded with fishing hook")
$pHookD     = New-Para("“Fishing hook, how interesting. Usually killers make it with knife. Probably he took a hook here, in the seaman house, but those shoe marks…”")
$pMessage   = New-Para("The message: “You are not a big fish!” or “Stay out of my business!” (second is stupid cuz he had already killed Klinton)")

$pPoliceman = New-Para("Policeman: It is not clear, right, Pier? Something is wrong here, I feel it. Take a number of his vital, take a minute to talk with him. When you are ready to tell who is a killer – call me!")

$pStage3    = New-Para("Stage 3.")

$pHint      = New-Para("Hint: Move to the office using a map, to talk with anybody whose number you have. *")

$pFreeMove  = New-Para("Now you can freely move to any location: Seaman Room, Pier and Office. ")

$pPier      = New-Para("At Pier you can find shoes in the sea and talk with a drunk-man Dude.")
$pOffice    = New-Para("In the office you can have a conversation with anybody whose phone number you have")
$pKlinton   = New-Para("In the Klinton’s house you can look and pieces of evidence.")

$pEmpty     = New-Para("")

$pLocPier   = New-Para("Location: Pier")

$pTabEnd    = New-Para("")

# ===== Now go back and apply paragraph-specific formatting / numbering =====

# "Items to touch" numbered list (decimal, 2 levels: "1." / "a.")
$pShoe.Style = "List Paragraph"
$pShoe.Range.ListFormat.ApplyNumberDefault()
$touchTemplate = $pShoe.Range.ListFormat.ListTemplate
$touchLevel2 = $touchTemplate.ListLevels.Item(2)
$touchLevel2.NumberStyle = 4
$touchLevel2.NumberFormat = "%2."

$pShoeDesc.Style = "List Paragraph"
$pShoeDesc.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pShoeDesc.Range.ListFormat.ListLevelNumber = 2

$pHarpoon.Style = "List Paragraph"
$pHarpoon.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pHarpoon.Range.ListFormat.ListLevelNumber = 1

$pHarpoonD.Style = "List Paragraph"
$pHarpoonD.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pHarpoonD.Range.ListFormat.ListLevelNumber = 2

$pHook.Style = "List Paragraph"
$pHook.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pHook.Range.ListFormat.ListLevelNumber = 1

$pHookD.Style = "List Paragraph"
$pHookD.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pHookD.Range.ListFormat.ListLevelNumber = 2

$pMessage.Style = "List Paragraph"
$pMessage.Range.ListFormat.ApplyListTemplateWithLevel($touchTemplate, $true, 1, $false, 0)
$pMessage.Range.ListFormat.ListLevelNumber = 2

# Bold "Stage 3."
$pStage3.Range.Font.Bold = $true

# "Hint" bold + underlined, rest of the sentence normal
$hintStart = $pHint.Range.Start
$hintWordRange = $d.Range($hintStart, $hintStart + 4)
$hintWordRange.Font.Bold = $true
$hintWordRange.Font.Underline = 1

# "In the office..." / "At Pier..." / "In the Klinton's..." numbered list (decimal, one level)
$pPier.Style = "List Paragraph"
$pPier.Range.ListFormat.ApplyNumberDefault()
$moveTemplate = $pPier.Range.ListFormat.ListTemplate

$pOffice.Style = "List Paragraph"
$pOffice.Range.ListFormat.ApplyListTemplateWithLevel($moveTemplate, $true, 1, $false, 0)

$pKlinton.Style = "List Paragraph"
$pKlinton.Range.ListFormat.ApplyListTemplateWithLevel($moveTemplate, $true, 1, $false, 0)

# Bold "Location: Pier"
$pLocPier.Range.Font.Bold = $true

# Final paragraph: bold tab, carries the _GoBack bookmark (moved from the "Items to touch" paragraph)
$tabRange = $pTabEnd.Range
$tabRange.Collapse(0)
$tabRange.InsertBefore("`t")
$pTabEnd.Range.Font.Bold = $true

$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bm.Delete()
$pTabEnd.Range.Collapse(0)
$d.Bookmarks.Add("_GoBack", $pTabEnd.Range)
